$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was updated
# from 45192 (2023-09-23) to 45202 (2023-10-03) for every data row
# (rows 2 through 173).
$lastRow = 173
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}
